# Generate Report for Handoff
# Updates the localization-status report with the newly generated
# handoff package (new GUID-named source file and new xliff hashes),
# along with the refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "2c19dca3-1354-459c-8ce3-f3904d0c48ed"
$newGuid = "82a0a388-0d8e-455e-9b26-ad91f2b58d3c"

$oldZhCnXlf = "$oldGuid.d56d0ecf6ea45a9698526a792fa678d06e9d16b2.zh-cn.xlf"
$newZhCnXlf = "$newGuid.4c8ae8774d7a21fb69634e20ae7bd7a7108515a5.zh-cn.xlf"

$oldDeDeXlf = "$oldGuid.d56d0ecf6ea45a9698526a792fa678d06e9d16b2.de-de.xlf"
$newDeDeXlf = "$newGuid.4c8ae8774d7a21fb69634e20ae7bd7a7108515a5.de-de.xlf"

$newGenerateDate = "2016-09-03 15:02:28"
$newZhCnHandoffDate = "2016-09-03 15:02:18"

# The source-file hyperlink target (points at the commit holding the source
# file on GitHub) is unchanged by this edit - only the visible display text
# changes to reflect the regenerated GUID-named file.
$sourceFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1cb66554abf100c9c02b062ab9da523ff21bc09/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $sourceFileUrl, "", "", "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = $newGenerateDate

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $sourceFileUrl, "", "", "$newGuid.md")
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $sourceFileUrl, "", "", "$newGuid.md")
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newGenerateDate
